$wb = $excel.ActiveWorkbook

# --- Sprint1: update projected (estimated) points for sprint 1 ---
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Range("E2").Value = 10
$sprint1.Range("F2").Value = 35
$sprint1.Range("E7").Value = 10
$sprint1.Range("F7").Value = 25

# --- Burndown: insert a new leading "Sprint"/"Start" label column ---
$burndown = $wb.Worksheets.Item("Burndown")
$burndown.Columns.Item(1).Insert()
$burndown.Range("A1").Value = "Sprint"
$burndown.Range("A2").Value = "Start"

# Re-point the Burndown chart series to the shifted columns (B = dates, C = points)
$chartObj = $burndown.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(,Burndown!`$B`$2:`$B`$7,Burndown!`$C`$2:`$C`$7,1)"
